$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 266.86667
$ws.Range("I33").Value = 228.78572
$ws.Range("K33").Value = 228.78572
$ws.Range("M33").Value = 0.2142800000000022
$ws.Range("H70").Value = 1700
$ws.Range("I70").Value = 1700
$ws.Range("K70").Value = 5100
$ws.Range("M70").Value = -4830
$ws.Range("H73").Value = 1700
$ws.Range("I73").Value = 1700
$ws.Range("K73").Value = 5100
$ws.Range("M73").Value = -4164
$ws.Range("H100").Value = 2147.6667
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459
$ws.Range("H107").Value = 504
$ws.Range("I107").Value = 170
$ws.Range("K107").Value = 170
$ws.Range("M107").Value = 1750
$ws.Range("H132").Value = 25555
$ws.Range("I132").Value = 25555
$ws.Range("K132").Value = 76665
$ws.Range("M132").Value = -74135
$ws.Range("H137").Value = 2743.1
$ws.Range("I137").Value = 1414.5
$ws.Range("J137").Value = 4736
$ws.Range("K137").Value = 4243.5
$ws.Range("L137").Value = 14208
$ws.Range("M137").Value = -1693.5
$ws.Range("N137").Value = -19308

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 320
$ws.Range("I2").Value = 320
$ws.Range("K2").Value = 320
$ws.Range("M2").Value = -207
$ws.Range("H38").Value = 2950
$ws.Range("I38").Value = 3000
$ws.Range("J38").Value = 2900
$ws.Range("K38").Value = 3000
$ws.Range("L38").Value = 2900
$ws.Range("M38").Value = -2533
$ws.Range("N38").Value = -3834
$ws.Range("H45").Value = 1398
$ws.Range("I45").Value = 1382.3334
$ws.Range("J45").Value = 1445
$ws.Range("K45").Value = 1382.3334
$ws.Range("L45").Value = 1445
$ws.Range("M45").Value = -1005.3334
$ws.Range("N45").Value = -2199
$ws.Range("H63").Value = 2273.4
$ws.Range("I63").Value = 1561
$ws.Range("J63").Value = 2748.3333
$ws.Range("K63").Value = 1561
$ws.Range("L63").Value = 2748.3333
$ws.Range("M63").Value = -875
$ws.Range("N63").Value = -4120.3333
$ws.Range("H66").Value = 2273.4
$ws.Range("I66").Value = 1561
$ws.Range("J66").Value = 2748.3333
$ws.Range("K66").Value = 7805
$ws.Range("L66").Value = 13741.6665
$ws.Range("M66").Value = -4373
$ws.Range("N66").Value = -20605.6665
$ws.Range("H102").Value = 1887.5
$ws.Range("I102").Value = 1887.5
$ws.Range("K102").Value = 1887.5
$ws.Range("M102").Value = -265.5
$ws.Range("H110").Value = 1151
$ws.Range("I110").Value = 1132.5
$ws.Range("K110").Value = 1132.5
$ws.Range("M110").Value = 912.5
$ws.Range("H116").Value = 320
$ws.Range("I116").Value = 320
$ws.Range("K116").Value = 320
$ws.Range("M116").Value = 1974
$ws.Range("H122").Value = 4200
$ws.Range("I122").Value = 1900
$ws.Range("K122").Value = 5700
$ws.Range("M122").Value = -3250

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 320
$ws.Range("I3").Value = 320
$ws.Range("K3").Value = 320
$ws.Range("M3").Value = -206
$ws.Range("H80").Value = 618.375
$ws.Range("J80").Value = 692.1667
$ws.Range("L80").Value = 692.1667
$ws.Range("N80").Value = -2688.1667
$ws.Range("H83").Value = 618.375
$ws.Range("J83").Value = 692.1667
$ws.Range("L83").Value = 3460.8335
$ws.Range("N83").Value = -13444.8335
$ws.Range("H86").Value = 1631.3334
$ws.Range("I86").Value = 1757.5
$ws.Range("J86").Value = 1379
$ws.Range("K86").Value = 1757.5
$ws.Range("L86").Value = 1379
$ws.Range("M86").Value = -634.5
$ws.Range("N86").Value = -3625
$ws.Range("H89").Value = 1631.3334
$ws.Range("I89").Value = 1757.5
$ws.Range("J89").Value = 1379
$ws.Range("K89").Value = 8787.5
$ws.Range("L89").Value = 6895
$ws.Range("M89").Value = -3171.5
$ws.Range("N89").Value = -18127
$ws.Range("H99").Value = 1074.4286
$ws.Range("I99").Value = 1074.4286
$ws.Range("K99").Value = 1074.4286
$ws.Range("M99").Value = 423.5714
$ws.Range("H107").Value = 1146.6666
$ws.Range("I107").Value = 1253.5
$ws.Range("J107").Value = 933
$ws.Range("K107").Value = 1253.5
$ws.Range("L107").Value = 933
$ws.Range("M107").Value = 666.5
$ws.Range("N107").Value = -4773

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2018.3334
$ws.Range("I16").Value = 2027.5
$ws.Range("K16").Value = 2027.5
$ws.Range("M16").Value = -1740.5
$ws.Range("H31").Value = 2984.625
$ws.Range("I31").Value = 1258.25
$ws.Range("J31").Value = 4711
$ws.Range("K31").Value = 1258.25
$ws.Range("L31").Value = 4711
$ws.Range("M31").Value = -963.25
$ws.Range("N31").Value = -5301
$ws.Range("H32").Value = 2575
$ws.Range("I32").Value = 100
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 100
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 216
$ws.Range("N32").Value = -10632
$ws.Range("H34").Value = 2984.625
$ws.Range("I34").Value = 1258.25
$ws.Range("J34").Value = 4711
$ws.Range("K34").Value = 1258.25
$ws.Range("L34").Value = 4711
$ws.Range("M34").Value = -1056.25
$ws.Range("N34").Value = -5115
$ws.Range("H62").Value = 4519.8
$ws.Range("I62").Value = 4374.75
$ws.Range("K62").Value = 4374.75
$ws.Range("M62").Value = -3750.75
$ws.Range("H65").Value = 4519.8
$ws.Range("I65").Value = 4374.75
$ws.Range("K65").Value = 21873.75
$ws.Range("M65").Value = -18753.75
$ws.Range("H75").Value = 60000
$ws.Range("J75").Value = 60000
$ws.Range("L75").Value = 60000
$ws.Range("N75").Value = -61996
$ws.Range("H78").Value = 60000
$ws.Range("J78").Value = 60000
$ws.Range("L78").Value = 180000
$ws.Range("N78").Value = -189984
$ws.Range("H105").Value = 805.5
$ws.Range("I105").Value = 500
$ws.Range("J105").Value = 1111
$ws.Range("K105").Value = 500
$ws.Range("L105").Value = 1111
$ws.Range("M105").Value = 1247
$ws.Range("N105").Value = -4605
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
$ws.Range("H113").Value = 2018.3334
$ws.Range("I113").Value = 2027.5
$ws.Range("K113").Value = 2027.5
$ws.Range("M113").Value = 142.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 259.66666
$ws.Range("J97").Value = 108
$ws.Range("L97").Value = 324
$ws.Range("N97").Value = -1316
$ws.Range("H108").Value = 500
$ws.Range("I108").Value = 500
$ws.Range("K108").Value = 1500
$ws.Range("M108").Value = 1380
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3037.25
$ws.Range("I132").Value = 2199.6667
$ws.Range("K132").Value = 6599.000100000001
$ws.Range("M132").Value = -4069.000100000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3474.8333
$ws.Range("I7").Value = 1500
$ws.Range("K7").Value = 1500
$ws.Range("M7").Value = -1388
$ws.Range("H32").Value = 3006.5
$ws.Range("I32").Value = 3006.5
$ws.Range("K32").Value = 3006.5
$ws.Range("M32").Value = -2689.5
$ws.Range("H61").Value = 2500
$ws.Range("I61").Value = 2500
$ws.Range("K61").Value = 2500
$ws.Range("M61").Value = -2298
$ws.Range("H93").Value = 4250
$ws.Range("I93").Value = 4000
$ws.Range("J93").Value = 4500
$ws.Range("K93").Value = 4000
$ws.Range("L93").Value = 4500
$ws.Range("M93").Value = -2752
$ws.Range("N93").Value = -6996
$ws.Range("H100").Value = 1420.4445
$ws.Range("I100").Value = 1420.4445
$ws.Range("K100").Value = 1420.4445
$ws.Range("M100").Value = -879.4445000000001
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = -330
$ws.Range("H126").Value = 3474.8333
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030
$ws.Range("H132").Value = 5545.2
$ws.Range("I132").Value = 5583
$ws.Range("J132").Value = 5394
$ws.Range("K132").Value = 16749
$ws.Range("L132").Value = 16182
$ws.Range("M132").Value = -14219
$ws.Range("N132").Value = -21242

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 993.25
$ws.Range("I113").Value = 993.25
$ws.Range("K113").Value = 2979.75
$ws.Range("M113").Value = -809.75
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("H126").Value = 1886.7142
$ws.Range("I126").Value = 1461.4
$ws.Range("K126").Value = 4384.200000000001
$ws.Range("M126").Value = -1914.200000000001
$ws.Range("N119").ClearContents()

Write-Output "applied $(236) writes and $(2) deletes"